# Regenerate the handoff status report.
#
# A fresh handoff pass completed for the files that were either
# "Ready for handoff" or had a failed handback transform. Their
# "Latest Handoff Date" / "Latest Handoff Datetime" timestamps move
# forward to the new handoff run; files still "In Translation" (or
# already handed back / in-sync) are untouched.

$wb = $excel.ActiveWorkbook

$rowsToRefresh = @(7, 10, 11, 12, 13, 14, 15, 16)

# --- Overview sheet: column D = "Latest Handoff Date" ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rowsToRefresh) {
    $overview.Cells.Item($r, 4).Value = "2016-24-19 06:24:21"
}

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rowsToRefresh) {
    $zhcn.Cells.Item($r, 5).Value = "2016-03-19 06:24:15"
}

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rowsToRefresh) {
    $dede.Cells.Item($r, 5).Value = "2016-03-19 06:24:21"
}
